$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 28-29, pushing existing rows 28..129 down to 30..131
$ws.Range("A28:A29").EntireRow.Insert()

# New row 28
$ws.Range("A28").Value = 3
$ws.Range("B28").Value = "Femacal de La Calera"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44707
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107011
$ws.Range("J28").Value = "Tuna"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 75
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 18000
$ws.Range("P28").Value = 18000
$ws.Range("Q28").Value = "$/caja 20 kilos"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 900
$ws.Range("T28").Value = 20

# New row 29
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44707
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107011
$ws.Range("J29").Value = "Tuna"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 70
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 16000
$ws.Range("Q29").Value = "$/caja 20 kilos"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 800
$ws.Range("T29").Value = 20

Write-Host "Inserted rows 28-29 and populated new data"
